$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(61, 1).Value = "2025-04-29 07:58:39"
$ws.Cells.Item(61, 2).Value = 182
